$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = 6260
$ws.Range("D21").Value = 5684562
$ws.Range("E21").Value = 908.0769968051118
$ws.Range("F21").Value = 8.661690678701618
$ws.Range("H21").Value = 29.75572544175129
